$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# MV data update: append the two new daily auction rows (02-11-2021 and 03-11-2021)
# A scratch cell far outside the used range is used to build the "Serie" labels as
# literal text (via a quoted formula) and then copy/paste-special the value into
# place, so the date-like strings are stored as plain text instead of being
# auto-converted into date serial numbers by Excel's normal text entry parsing.
$scratch = $ws.Cells.Item(1000, 1)

$scratch.Formula = "=""02-11-2021"""
$scratch.Copy()
$ws.Cells.Item(142, 1).PasteSpecial(-4163)
$scratch.ClearContents()

$ws.Cells.Item(142, 2).Value = 1000000
$ws.Cells.Item(142, 3).Value = 1815000
$ws.Cells.Item(142, 4).Value = 1000000
$ws.Cells.Item(142, 5).Value = 615000
$ws.Cells.Item(142, 6).Value = 385000
$ws.Cells.Item(142, 7).Value = 2.75

$scratch.Formula = "=""03-11-2021"""
$scratch.Copy()
$ws.Cells.Item(143, 1).PasteSpecial(-4163)
$scratch.ClearContents()

$ws.Cells.Item(143, 2).Value = 1200000
$ws.Cells.Item(143, 3).Value = 1585000
$ws.Cells.Item(143, 4).Value = 1200000
$ws.Cells.Item(143, 5).Value = 875000
$ws.Cells.Item(143, 6).Value = 325000
$ws.Cells.Item(143, 7).Value = 2.75
